$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (shared-string content renames) ---
# "Período" field: dt_data_base -> dt_base (used identically by E4 and F4)
$ws.Range("E4").Value = "dt_base"
$ws.Range("F4").Value = "dt_base"

# "Operações" field: qt_numero_de_operacoes -> qt_operacoes
$ws.Range("E14").Value = "qt_operacoes"

# --- F14 becomes a truly blank cell (was an empty-string / quotePrefix cell) ---
# Copy the plain blank-cell format from a neighbour and clear the value so
# the quotePrefix style (and the now-unused empty shared string) go away.
$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").ClearContents()
$excel.CutCopyMode = 0

# --- Normalize C25 / C26 style (drop the alternate "font2" border style) ---
$ws.Range("D25").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New field row: "Período (Ano)" / dt_ano ---
$ws.Range("A27").Value = "N/A"
$ws.Range("B27").Value = "Per" + [char]0x00ed + "odo (Ano)"
$ws.Range("F27").Value = "dt_ano"

# --- Row heights: rows 1-23 shrink from 19.5 to 18.75 ---
for ($i = 1; $i -le 23; $i++) {
    $ws.Rows.Item($i).RowHeight = 18.75
}
